$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header cell in A1 was renamed from "CarKilometers" to "Name"
$ws.Range("A1").Value = "Name"
